# Fruta / hortaliza, semanal
# Insert a new weekly record at row 103 (pushing the existing rows 103-163
# down to 104-164) and populate the new row with the latest price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 103:163 down to 104:164, leaving a blank row 103 behind.
$ws.Rows(103).Insert()

# Fill in the new row 103 with the new weekly entry.
$ws.Cells.Item(103, 1).Value = 4
$ws.Cells.Item(103, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(103, 3).Value = "Los Lagos"
$ws.Cells.Item(103, 4).Value = 44606
$ws.Cells.Item(103, 5).Value = 10
$ws.Cells.Item(103, 6).Value = "Fruta"
$ws.Cells.Item(103, 7).Value = 100108
$ws.Cells.Item(103, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(103, 9).Value = 100108002
$ws.Cells.Item(103, 10).Value = "Mango"
$ws.Cells.Item(103, 11).Value = "Sin especificar"
$ws.Cells.Item(103, 12).Value = "Primera"
$ws.Cells.Item(103, 13).Value = 80
$ws.Cells.Item(103, 14).Value = 7500
$ws.Cells.Item(103, 15).Value = 8000
$ws.Cells.Item(103, 16).Value = 7750
$ws.Cells.Item(103, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(103, 18).Value = "Perú"
$ws.Cells.Item(103, 19).Value = 1938
$ws.Cells.Item(103, 20).Value = 4
